$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the cells we will touch as Text so Excel keeps the
# exact literal strings (prices/percentages) instead of coercing
# them into floating point numbers.
$targets = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($ref in $targets) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated price / volume(1h) values.
$ws.Range("D2").Value = "325.26"
$ws.Range("E2").Value = "-2.20%"
$ws.Range("D3").Value = "44.50"
$ws.Range("E3").Value = "0.56%"
$ws.Range("D4").Value = "5.495"
$ws.Range("E4").Value = "-5.77%"
$ws.Range("D5").Value = "0.08052"
$ws.Range("E5").Value = "-3.44%"
$ws.Range("D6").Value = "8.657"
$ws.Range("E6").Value = "-1.81%"
$ws.Range("D7").Value = "1.913"
$ws.Range("E7").Value = "-3.11%"
$ws.Range("D8").Value = "4.285"
$ws.Range("E8").Value = "-4.88%"
$ws.Range("D9").Value = "2.715"
$ws.Range("E9").Value = "-6.33%"
$ws.Range("D10").Value = "0.9420"
$ws.Range("E10").Value = "0.66%"
$ws.Range("E11").Value = "-7.74%"
$ws.Range("D12").Value = "0.1870"
$ws.Range("E12").Value = "-4.00%"
$ws.Range("D13").Value = "0.09988"
$ws.Range("E13").Value = "4.28%"
$ws.Range("D14").Value = "0.04230"
$ws.Range("E14").Value = "5.82%"
$ws.Range("D15").Value = "0.1065"
$ws.Range("E15").Value = "-0.11%"
$ws.Range("D16").Value = "0.001280"
$ws.Range("D17").Value = "0.005891"
$ws.Range("E17").Value = "-2.53%"
$ws.Range("D18").Value = "3.588"
$ws.Range("E18").Value = "2.33%"
$ws.Range("D19").Value = "0.3485"
$ws.Range("E19").Value = "-0.72%"
$ws.Range("D20").Value = "8.523"
$ws.Range("E20").Value = "-5.01%"
$ws.Range("D21").Value = "0.1372"
$ws.Range("E21").Value = "-0.03%"
$ws.Range("D22").Value = "0.2529"
$ws.Range("E22").Value = "-1.66%"
$ws.Range("D23").Value = "0.04244"
$ws.Range("E23").Value = "-4.00%"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "-1.76%"
$ws.Range("D25").Value = "0.004555"
$ws.Range("E25").Value = "3.53%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.78%"
$ws.Range("D27").Value = "0.0003991"
$ws.Range("E27").Value = "-0.07%"
$ws.Range("D39").Value = "0.02638"
$ws.Range("E39").Value = "-5.63%"
$ws.Range("D40").Value = "0.05473"
$ws.Range("E40").Value = "-4.06%"
$ws.Range("D41").Value = "0.007697"
$ws.Range("E41").Value = "-2.65%"
$ws.Range("E42").Value = "-2.27%"
$ws.Range("D43").Value = "0.007179"
$ws.Range("E43").Value = "-20.21%"
$ws.Range("D44").Value = "0.002048"
$ws.Range("E44").Value = "-2.60%"
$ws.Range("D45").Value = "0.008663"
$ws.Range("E45").Value = "-17.45%"
$ws.Range("D46").Value = "0.00007108"
$ws.Range("E46").Value = "-2.18%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").Value = "0.003606"
$ws.Range("E48").Value = "11.15%"
$ws.Range("D49").Value = "0.002271"
$ws.Range("E49").Value = "-0.40%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.06%"

Write-Host "Updated $($targets.Count) cells."
